# Auto-applied update: dades i banners [2026-02-20 23:50]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 23:48:18"
$ws.Range("E3").Value = "2026-02-20 23:48:20"
$ws.Range("O3").Value = "-5.1 °C"
$ws.Range("E4").Value = "2026-02-20 23:48:23"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "62%"
$ws.Range("J4").Value = "1023.3 hPa"
$ws.Range("N4").Value = "3.1 °C 23:05 TU"
$ws.Range("O4").Value = "9.4 °C"
$ws.Range("E5").Value = "2026-02-20 23:48:26"
$ws.Range("E6").Value = "2026-02-20 23:48:28"
$ws.Range("J6").Value = "1023.2 hPa"
$ws.Range("E7").Value = "2026-02-20 23:48:31"
$ws.Range("J7").Value = "1023.1 hPa"
$ws.Range("E8").Value = "2026-02-20 23:48:33"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "61%"
$ws.Range("J8").Value = "1023.4 hPa"
$ws.Range("E9").Value = "2026-02-20 23:48:36"
$ws.Range("E10").Value = "2026-02-20 23:48:38"
$ws.Range("O10").Value = "7.2 °C"
$ws.Range("E11").Value = "2026-02-20 23:48:41"
$ws.Range("O11").Value = "9.1 °C"
$ws.Range("E12").Value = "2026-02-20 23:48:43"
$ws.Range("E13").Value = "2026-02-20 23:48:46"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "48%"
$ws.Range("J13").Value = "1024.8 hPa"
$ws.Range("N13").Value = "-2.2 °C 23:10 TU"
$ws.Range("O13").Value = "5.5 °C"
$ws.Range("E14").Value = "2026-02-20 23:48:48"
$ws.Range("O14").Value = "11.5 °C"
$ws.Range("E15").Value = "2026-02-20 23:48:50"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "47%"
$ws.Range("E16").Value = "2026-02-20 23:48:53"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "49%"
$ws.Range("O16").Value = "-2.9 °C"
$ws.Range("E17").Value = "2026-02-20 23:48:55"
$ws.Range("K17").Value = "5.1 MJ/m2"
$ws.Range("M17").Value = "6.7 °C 23:24 TU"
$ws.Range("E18").Value = "2026-02-20 23:48:58"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "80%"
$ws.Range("J18").Value = "1023.5 hPa"
$ws.Range("O18").Value = "7.4 °C"
$ws.Range("E19").Value = "2026-02-20 23:49:00"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "70%"
$ws.Range("E20").Value = "2026-02-20 23:49:03"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "56%"
$ws.Range("E21").Value = "2026-02-20 23:49:05"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "40%"
$ws.Range("J21").Value = "1023.6 hPa"
$ws.Range("N21").Value = "2.3 °C 23:03 TU"
$ws.Range("O21").Value = "8.7 °C"
$ws.Range("E22").Value = "2026-02-20 23:49:08"
$ws.Range("O22").Value = "-3.8 °C"
$ws.Range("E23").Value = "2026-02-20 23:49:10"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "64%"
$ws.Range("O23").Value = "-4.3 °C"
$ws.Range("E24").Value = "2026-02-20 23:49:13"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "69%"
$ws.Range("J24").Value = "1025.9 hPa"
$ws.Range("N24").Value = "4.2 °C 23:29 TU"
$ws.Range("O24").Value = "9.2 °C"
$ws.Range("E25").Value = "2026-02-20 23:49:15"
$ws.Range("E26").Value = "2026-02-20 23:49:18"
$ws.Range("J26").Value = "1022.3 hPa"
$ws.Range("E27").Value = "2026-02-20 23:49:20"
$ws.Range("E28").Value = "2026-02-20 23:49:22"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "68%"
$ws.Range("J28").Value = "1023.8 hPa"
$ws.Range("O28").Value = "6.6 °C"
$ws.Range("E29").Value = "2026-02-20 23:49:25"
$ws.Range("E30").Value = "2026-02-20 23:49:27"
$ws.Range("J30").Value = "1023.0 hPa"
$ws.Range("E31").Value = "2026-02-20 23:49:29"
$ws.Range("J31").Value = "1022.2 hPa"
$ws.Range("E32").Value = "2026-02-20 23:49:32"
$ws.Range("O32").Value = "4.3 °C"
$ws.Range("E33").Value = "2026-02-20 23:49:34"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "45%"
$ws.Range("J33").Value = "1024.0 hPa"
$ws.Range("N33").Value = "1.3 °C 23:02 TU"
$ws.Range("O33").Value = "5.6 °C"
$ws.Range("E34").Value = "2026-02-20 23:49:37"
$ws.Range("K34").Value = "10.4 MJ/m2"
$ws.Range("O34").Value = "1.3 °C"
$ws.Range("E35").Value = "2026-02-20 23:49:39"
$ws.Range("J35").Value = "1027.3 hPa"
$ws.Range("K35").Value = "10.8 MJ/m2"
$ws.Range("E36").Value = "2026-02-20 23:49:42"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "44%"
$ws.Range("J36").Value = "1023.1 hPa"
$ws.Range("O36").Value = "13.9 °C"
$ws.Range("E37").Value = "2026-02-20 23:49:44"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "70%"
$ws.Range("J37").Value = "1025.4 hPa"
$ws.Range("O37").Value = "4.2 °C"
$ws.Range("E38").Value = "2026-02-20 23:49:47"
$ws.Range("E39").Value = "2026-02-20 23:49:49"
$ws.Range("E40").Value = "2026-02-20 23:49:52"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "40%"
$ws.Range("J40").Value = "1024.3 hPa"
$ws.Range("O40").Value = "9.8 °C"
$ws.Range("E41").Value = "2026-02-20 23:49:54"
$ws.Range("J41").Value = "1023.6 hPa"
$ws.Range("E42").Value = "2026-02-20 23:49:57"
$ws.Range("E43").Value = "2026-02-20 23:49:59"
$ws.Range("K43").Value = "6.8 MJ/m2"
$ws.Range("E44").Value = "2026-02-20 23:50:01"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "73%"
$ws.Range("M44").Value = "0.4 °C 23:29 TU"
$ws.Range("O44").Value = "-4.1 °C"
$ws.Range("E45").Value = "2026-02-20 23:50:04"
$ws.Range("J45").Value = "1030.2 hPa"
$ws.Range("E46").Value = "2026-02-20 23:50:06"
$ws.Range("N46").Value = "7.1 °C 23:21 TU"
$ws.Range("O46").Value = "11.6 °C"
